# Apply updated cryptos list values (price + volume change %) per row
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cell, $text) {
    $origStyle = $cell.Style
    $cell.NumberFormat = "@"
    $cell.Value = $text
    $cell.Style = $origStyle
}

Set-TextValue $ws.Range("D2") "63.560.23"
Set-TextValue $ws.Range("E2") "  +0.35%  "
Set-TextValue $ws.Range("D3") "2.583.88"
Set-TextValue $ws.Range("E3") "  -0.95%  "
Set-TextValue $ws.Range("D4") "1.00"
Set-TextValue $ws.Range("E4") "  +0.14%  "
Set-TextValue $ws.Range("D5") "588.12"
Set-TextValue $ws.Range("E5") "  -0.24%  "
Set-TextValue $ws.Range("D6") "144.80"
Set-TextValue $ws.Range("E6") "  -3.16%  "
Set-TextValue $ws.Range("D7") "1.00"
Set-TextValue $ws.Range("E7") "  +0.12%  "
Set-TextValue $ws.Range("E8") "  -1.62%  "
Set-TextValue $ws.Range("E9") "  -3.46%  "
Set-TextValue $ws.Range("D10") "5.62"
Set-TextValue $ws.Range("E10") "  -0.98%  "
Set-TextValue $ws.Range("E11") "  -0.24%  "
Set-TextValue $ws.Range("D12") "0.352"
Set-TextValue $ws.Range("E12") "  -1.82%  "
Set-TextValue $ws.Range("D13") "27.38"
Set-TextValue $ws.Range("E13") "  -1.81%  "
Set-TextValue $ws.Range("D14") "3.048.54"
Set-TextValue $ws.Range("E14") "  -0.78%  "
Set-TextValue $ws.Range("D15") "63.445.10"
Set-TextValue $ws.Range("E15") "  +0.18%  "
Set-TextValue $ws.Range("E16") "  -2.69%  "
Set-TextValue $ws.Range("D17") "2.579.62"
Set-TextValue $ws.Range("E17") "  -0.45%  "
Set-TextValue $ws.Range("D18") "11.13"
Set-TextValue $ws.Range("E18") "  -3.16%  "
Set-TextValue $ws.Range("D19") "343.34"
Set-TextValue $ws.Range("E19") "  -0.83%  "
Set-TextValue $ws.Range("D20") "4.32"
Set-TextValue $ws.Range("E20") "  -3.32%  "
Set-TextValue $ws.Range("D21") "6.64"
Set-TextValue $ws.Range("E21") "  -3.83%  "
Set-TextValue $ws.Range("E22") "  +0.22%  "
Set-TextValue $ws.Range("D23") "68.41"
Set-TextValue $ws.Range("E23") "  +2.21%  "
Set-TextValue $ws.Range("D24") "1.57"
Set-TextValue $ws.Range("E24") "  +5.43%  "
Set-TextValue $ws.Range("D25") "1.61"
Set-TextValue $ws.Range("E25") "  -1.62%  "
Set-TextValue $ws.Range("E26") "  -3.96%  "
Set-TextValue $ws.Range("E27") "  +0.12%  "
Set-TextValue $ws.Range("D28") "7.94"
Set-TextValue $ws.Range("E28") "  -3.66%  "
Set-TextValue $ws.Range("D29") "8.26"
Set-TextValue $ws.Range("E29") "  -3.46%  "
Set-TextValue $ws.Range("E30") "  -2.90%  "
Set-TextValue $ws.Range("D31") "471.75"
Set-TextValue $ws.Range("E31") "  +0.58%  "
Set-TextValue $ws.Range("D32") "0.0₃0802"
Set-TextValue $ws.Range("E32") "  -4.27%  "
Set-TextValue $ws.Range("D33") "1.68"
Set-TextValue $ws.Range("E33") "  +1.67%  "
Set-TextValue $ws.Range("D34") "176.65"
Set-TextValue $ws.Range("E34") "  -0.32%  "
Set-TextValue $ws.Range("E35") "  +0.28%  "
Set-TextValue $ws.Range("E36") "  -2.05%  "
Set-TextValue $ws.Range("D37") "18.93"
Set-TextValue $ws.Range("E37") "  -2.33%  "
Set-TextValue $ws.Range("D38") "4.52"
Set-TextValue $ws.Range("E38") "  -3.33%  "
Set-TextValue $ws.Range("E39") "  +0.01%  "
Set-TextValue $ws.Range("E40") "  -1.94%  "
Set-TextValue $ws.Range("D41") "162.05"
Set-TextValue $ws.Range("E41") "  +5.52%  "
Set-TextValue $ws.Range("D42") "40.06"
Set-TextValue $ws.Range("E42") "  +0.98%  "
Set-TextValue $ws.Range("D43") "3.73"
Set-TextValue $ws.Range("E43") "  -3.77%  "
Set-TextValue $ws.Range("D44") "21.89"
Set-TextValue $ws.Range("E44") "  +3.07%  "
Set-TextValue $ws.Range("D45") "0.632"
Set-TextValue $ws.Range("E45") "  +2.13%  "
Set-TextValue $ws.Range("D46") "0.0538"
Set-TextValue $ws.Range("E46") "  -3.51%  "
Set-TextValue $ws.Range("D47") "0.0962"
Set-TextValue $ws.Range("E47") "  -1.80%  "
Set-TextValue $ws.Range("D48") "0.0238"
Set-TextValue $ws.Range("E48") "  -2.09%  "
Set-TextValue $ws.Range("D49") "18.26"
Set-TextValue $ws.Range("E49") "  -2.73%  "
Set-TextValue $ws.Range("D50") "11.35"
Set-TextValue $ws.Range("E50") "  -0.40%  "
Set-TextValue $ws.Range("E51") "  -4.58%  "
